$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy style (bold, border, centered) from A3 to A4 so it matches existing date cells
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122) # xlPasteFormats

# Write "07/02/2022" as a formula producing a text string, then collapse it down to
# a plain value via copy / paste-values. This keeps the cell's stored type as text
# (matching the shared string used by A2/A3) instead of Excel auto-parsing the
# date-like text into a date serial number when typed directly into .Value.
$ws.Range("A4").Formula = '="07/02/2022"'
$ws.Range("A4").Copy()
$ws.Range("A4").PasteSpecial(-4163) # xlPasteValues

# Numeric values for row 4
$ws.Range("B4").Value = 3805.269454
$ws.Range("C4").Value = 11335.950417
$ws.Range("D4").Value = 0.55
$ws.Range("E4").Value = 10.65
